$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# SCHEMATIC section: mark the "RF module needs 3.3V" row as 100% done.
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = "0%"

# New row describing the part verification task needed for the schematic.
$ws.Range("B19").Value = "Part Verification"
$ws.Range("C19").Value = "Size, Value, Make, Female"
$ws.Rows("19:19").RowHeight = 30

# Update the view state: scroll so row 4 is at the top and select C19.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C19").Select()
